# Slide 4, Title shape ("Title 1"): rewrite the title paragraph text.
# Before: "ส่วนประกอบของ "(th-TH) + "Class"(en-US)
# After : " Class "(en-US) + "ประกอบด้วยสมาชิก 2 ประเภท"(th-TH)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)
$sh = $s.Shapes.Item(1)
$tr = $sh.TextFrame.TextRange

# Replace the whole paragraph with the first (English) run's text.
$tr.Text = " Class "

# Append the second (Thai) run's text; InsertAfter gives back a range over
# just the newly-inserted text, which we use to give it its own run/lang.
$thaiRange = $tr.InsertAfter("ประกอบด้วยสมาชิก 2 ประเภท")
$thaiRange.LanguageID = "th-TH"

# Set the language of the first run last, so it doesn't get clobbered by
# the split above.
$tr.LanguageID = "en-US"
